# Refresh the cryptos list: update each coin's Price (col D) and
# Volume(1h) (col E) text, and (for rows 42/43) swap Fetch.AI/Maker.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. "69.709.67", "1.00").
# A plain .Value assignment lets Excel auto-coerce these to real numbers
# (dropping formatting like trailing zeros or the thousands dots), so we
# force the cell to Text first, then reset the style back to "Normal" so
# no stray text-format style index is left attached to the cell.
function Set-TextValue {
    param($Cell, $Text)
    $range = $ws.Range($Cell)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = "Normal"
}

Set-TextValue 'D2' '69.709.67'
$ws.Range('E2').Value = '  -0.68%  '
Set-TextValue 'D3' '3.561.31'
$ws.Range('E3').Value = '  -1.10%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue 'D5' '576.02'
$ws.Range('E5').Value = '  -3.21%  '
Set-TextValue 'D6' '187.48'
$ws.Range('E6').Value = '  -2.98%  '
$ws.Range('E7').Value = '  -3.36%  '
Set-TextValue 'D8' '3.556.71'
$ws.Range('E8').Value = '  -1.07%  '
Set-TextValue 'D9' '1.00'
$ws.Range('E9').Value = '  +0.02%  '
Set-TextValue 'D10' '0.176'
$ws.Range('E10').Value = '  -3.03%  '
$ws.Range('E11').Value = '  -1.20%  '
Set-TextValue 'D12' '55.78'
$ws.Range('E12').Value = '  -4.18%  '
Set-TextValue 'D13' '0.0000297'
$ws.Range('E13').Value = '  +0.48%  '
Set-TextValue 'D14' '9.63'
$ws.Range('E14').Value = '  -1.40%  '
Set-TextValue 'D15' '4.129.52'
$ws.Range('E15').Value = '  -0.98%  '
Set-TextValue 'D16' '19.84'
$ws.Range('E16').Value = '  +2.53%  '
Set-TextValue 'D17' '3.552.75'
$ws.Range('E17').Value = '  -1.17%  '
Set-TextValue 'D18' '69.521.70'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('E19').Value = '  -1.09%  '
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('E21').Value = '  -1.96%  '
Set-TextValue 'D22' '473.03'
$ws.Range('E22').Value = '  -5.36%  '
Set-TextValue 'D23' '19.35'
$ws.Range('E23').Value = '  +13.82%  '
Set-TextValue 'D24' '5.03'
$ws.Range('E24').Value = '  -8.91%  '
$ws.Range('E25').Value = '  -3.48%  '
Set-TextValue 'D26' '88.18'
$ws.Range('E26').Value = '  -3.44%  '
Set-TextValue 'D27' '3.03'
$ws.Range('E27').Value = '  -2.28%  '
Set-TextValue 'D28' '10.92'
$ws.Range('E28').Value = '  -2.64%  '
Set-TextValue 'D29' '9.32'
$ws.Range('E29').Value = '  -0.52%  '
Set-TextValue 'D30' '31.94'
$ws.Range('E30').Value = '  -1.39%  '
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('E32').Value = '  +2.04%  '
$ws.Range('E33').Value = '  -1.55%  '
Set-TextValue 'D34' '65.62'
$ws.Range('E34').Value = '  +0.30%  '
Set-TextValue 'D35' '573.28'
$ws.Range('E35').Value = '  -7.02%  '
Set-TextValue 'D36' '38.59'
$ws.Range('E36').Value = '  +1.36%  '
$ws.Range('E37').Value = '  -0.03%  '
Set-TextValue 'D38' '0.0₃0795'
$ws.Range('E38').Value = '  -4.85%  '
$ws.Range('E39').Value = '  -1.59%  '
$ws.Range('E40').Value = '  -5.90%  '
$ws.Range('E41').Value = '  -5.89%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D42' '3.209.20'
$ws.Range('E42').Value = '  -3.89%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D43' '2.84'
$ws.Range('E43').Value = '  +5.54%  '
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('E45').Value = '  +10.34%  '
$ws.Range('E46').Value = '  -0.78%  '
Set-TextValue 'D47' '9.43'
$ws.Range('E47').Value = '  +3.98%  '
Set-TextValue 'D48' '3.32'
$ws.Range('E48').Value = '  +1.60%  '
$ws.Range('E49').Value = '  -0.81%  '
Set-TextValue 'D50' '0.998'
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('E51').Value = '  -3.60%  '
